$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.99999999456735889,
    0.99711216730554431,
    0.98948684640679763,
    0.99166734508190224,
    0.98067295885175465,
    0.95441912211914426,
    0.94909685203559357,
    0.94225898402071739,
    0.93426789092828799,
    0.92719401828717685,
    0.9262161835382261,
    0.92343628294861546,
    0.91214893169967071,
    0.90798150636604902,
    0.90538998059014575,
    0.90288346086217008,
    0.89917550334666108,
    0.89806658583142585,
    0.99605895384219301,
    0.98894192372897183,
    0.98754343478820483,
    0.9862789312671395,
    0.97981428552881278,
    0.96679356346228573,
    0.96033659177583164,
    0.95177043525141802,
    0.950204121346655,
    0.94372687805836852,
    0.93963409622685612,
    0.93874450114130192,
    0.94531299964907267,
    0.94363372311069582,
    0.94311371296800006
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
